# Rename the inline picture shapes that live in the document's footers /
# first-page header. Only the `name` metadata ("image2.png" -> "image1.png"
# for the two Pearson logo instances in the footers, and "image1.jpg" ->
# "image2.jpg" for the BTEC logo in the first-page header) changes; the
# pictures themselves, their ids, sizes, and embedded media stay the same.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# wdHeaderFooterPrimary = 1 (default), wdHeaderFooterFirstPage = 2 (first page)

# Default footer -> Pearson logo (docPr id="1"): image2.png -> image1.png
$ftrDefault = $sec.Footers.Item(1)
$shpFtrDefault = $ftrDefault.Range.InlineShapes.Item(1)
$shpFtrDefault.Name = "image1.png"

# First-page footer -> Pearson logo (docPr id="2"): image2.png -> image1.png
$ftrFirst = $sec.Footers.Item(2)
$shpFtrFirst = $ftrFirst.Range.InlineShapes.Item(1)
$shpFtrFirst.Name = "image1.png"

# First-page header -> BTEC logo (docPr id="3"): image1.jpg -> image2.jpg
$hdrFirst = $sec.Headers.Item(2)
$shpHdrFirst = $hdrFirst.Range.InlineShapes.Item(1)
$shpHdrFirst.Name = "image2.jpg"

Write-Output ("Footer default image name: " + $shpFtrDefault.Name)
Write-Output ("Footer first-page image name: " + $shpFtrFirst.Name)
Write-Output ("Header first-page image name: " + $shpHdrFirst.Name)
